# dbo_counmet.xlsx was re-saved by a newer Excel build. Most of the
# resulting diff is Excel-internal bookkeeping (fileVersion/rupBuild,
# xr:revisionPtr/documentId, mc:AlternateContent author path, the
# "Office 2007 - 2010" -> "Office" theme upgrade, x14ac dyDescent/
# knownFonts markers, customXml part renumbering, ...) that isn't
# exposed through the Excel object model - it's stamped by the app
# itself when it writes the file, not something a user/script action
# produces. The two user-visible, OM-reachable edits are reproduced
# below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's default/"Normal" font changed from "MS Sans Serif" to
# "Calibri" - update the workbook's Normal style font so every cell
# (none of which carry an explicit font override) renders in Calibri,
# matching the visible effect of the diff's font rename.
$wb.Styles.Item("Normal").Font.Name = "Calibri"

# The saved cursor position moved from B1 to B12.
[void]$ws.Range("B12").Select()
